# Apply the LeadProcessWF_UsingSequence.xlsx maintenance edit:
#  - D2: "mir1" -> "mir2"
#  - E2: replace the nucleotide sequence text with the new (trimmed) sequence
#  - F2: "Butler KM" -> "Bharitkar S, Mendel"
#  - Active selection on Sheet1 moves from G2 to D2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Match the original authoring order (shared-string table append order):
# F2, then E2, then D2.
$ws.Range("F2").Value = "Bharitkar S, Mendel"

$newSequence = "GGAGGAGCGTGCGCGGGGGCCCCGGGA`nGACGGCGGCGGTGGCGGCGCGGGCAGAGCAAGGACGCGGCGGATCCCACTCGCACAGCAGCGCACTCGGT`nGCCCCGCGCAGGGTCGCGATGCTGCCCGGTTTGGCACTGCTCCTGCTGGCCGCCTGGACGGCTCGGGCGC`nTGGAGGTGGGTGCCGCGCCTCGGAAGGCGGGGGGAGGCTGCACGGTGGGGACGCGATACCCCCCAAGACC`nTTAACCCAAGTCTTTAATGCAGAGAAGCCGGGGGTCCGTCAATGGGACCCCTCTCCTCTCCGCCCCCGCT`nTGCGGACGTCCAGCGCATCCCCGCTTTCGGCCCAGCCCTGCCCCAGGGAGTCGCGCTCCGGCCCGCTGAG`nAGGGAGCGGGCGAGGCGCTGGTCTCCCTGGTTCCGCGCCAGCCCGGGGCGAGAAGGGTAGGGGGCGACCC`nTGAGCCCAGACCCCGACTTAGTCCCTGCCTTGGAAGCGGGGGTCGGGGGAGGCGAGAGACATTCAGACAG"
$ws.Range("E2").Value = $newSequence

$ws.Range("D2").Value = "mir2"

# Writing the longer/shorter sequence text triggers the engine's automatic
# row-autofit; restore the original authored row height (wrap-text row 2).
$ws.Rows.Item(2).RowHeight = 40.5

$ws.Range("D2").Select()
